# TIMES_DH_Demand_2050.xlsx - "Correction on heat sector data (1)"
#
# The old sheet had a spurious first row (B1 = "District heat"), then a
# header row (A2 = "Region", B2 = "District heating") followed by the
# per-country district-heat demand figures in rows 3:33.
#
# The fix removes that spurious first row (shifting the header up to row 1
# and the country rows up to 2:32) and corrects a number of the per-country
# demand values. The B column also gets a tighter number format (one
# decimal place) and is resized to fit the now much shorter numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the stray first row; everything below shifts up by one.
$ws.Rows.Item(1).Delete()

# The "District heating" header (now B1) loses the special formatting it
# inherited from the old B2 cell it shifted up from.
$ws.Range("B1").ClearFormats()

# Corrected per-country "District heating" demand values (column B),
# now sitting in rows 2:32 after the row shift above.
$ws.Range("B2").Value = 58.931000000000004
$ws.Range("B3").Value = 17.943000000000001
$ws.Range("B4").Value = 16.956
$ws.Range("B5").Value = 13.7
$ws.Range("B6").Value = 2.23
$ws.Range("B7").Value = 95.849000000000004
$ws.Range("B8").Value = 473.41200000000003
$ws.Range("B9").Value = 42.966999999999992
$ws.Range("B10").Value = 17.175000000000001
$ws.Range("B11").Value = 0.17499999999999999
$ws.Range("B12").Value = 7.0380000000000003
$ws.Range("B13").Value = 48.101999999999997
$ws.Range("B14").Value = 24.015000000000001
$ws.Range("B15").Value = 9.1750000000000007
$ws.Range("B16").Value = 35.707999999999998
$ws.Range("B17").Value = 0.92100000000000004
$ws.Range("B18").Value = 9.9589999999999996
$ws.Range("B19").Value = 249.614
$ws.Range("B20").Value = 29.657
$ws.Range("B21").Value = 0.93800000000000006
$ws.Range("B22").Value = 21.543000000000003
$ws.Range("B23").Value = 0.31900000000000001
$ws.Range("B24").Value = 63.656000000000006
$ws.Range("B25").Value = 25.529
$ws.Range("B26").Value = 189.928
$ws.Range("B27").Value = 13.797000000000001
$ws.Range("B28").Value = 108.97500000000001
$ws.Range("B29").Value = 110.334
$ws.Range("B30").Value = 6.6379999999999999
$ws.Range("B31").Value = 24.414000000000001
$ws.Range("B32").Value = 21.408999999999999

# Tighter display format (1 decimal place instead of 0) for the values.
$ws.Range("B2:B32").NumberFormat = "#,##0.0"

# Column B no longer needs to fit "District heating" widths of data -
# shrink it to fit the now-short numeric values.
$ws.Columns.Item(2).ColumnWidth = 4

# Matches the saved selection left behind in the authored workbook.
$ws.Range("D4").Select()
